$p = $ppt.ActivePresentation

# --- Slide 3: "Content Placeholder 2" -----------------------------------
# Merge the two runs "npm run " + "coverage" into a single run
# "npm run coverage" (same visible text, diff only changes the XML run
# structure).
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(3)
$tr3 = $shp3.TextFrame.TextRange

$needle3 = "npm run coverage"
$full3 = $tr3.Text
$idx3 = $full3.IndexOf($needle3)
if ($idx3 -ge 0) {
    $sub3 = $tr3.Characters($idx3 + 1, $needle3.Length)
    $sub3.Text = $needle3
}

# --- Slide 5: "Content Placeholder 2" -----------------------------------
# "differently than Angular 1.x" -> "manually"
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange

$needle5 = "differently than Angular 1.x"
$full5 = $tr5.Text
$idx5 = $full5.IndexOf($needle5)
if ($idx5 -ge 0) {
    $sub5 = $tr5.Characters($idx5 + 1, $needle5.Length)
    $sub5.Text = "manually"
}
